$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "64.607.91"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.368.16"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "563.32"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "176.62"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.623"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.08%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.359.05"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E10").Value = "  +8.46%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.632"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.50%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "55.31"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000276"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("E14").Value = "  +1.49%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.904.65"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "18.34"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.365.69"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  +1.23%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "64.548.67"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.59%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.990"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "462.54"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +13.29%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.83"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +9.86%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "86.30"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.53%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "13.53"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  +3.81%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.82"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "30.18"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.74"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.30%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "11.51"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "579.70"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").Value = "  +1.53%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "59.67"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  -7.24%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "35.95"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0₃0759"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("E40").Value = "  +1.98%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.097.38"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("E43").Value = "  +0.03%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.84"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.54%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.51"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  +2.12%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +1.70%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "135.97"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.07%  "
